# Apply crypto price/volume updates per commit "Updated cryptos list on Sun Oct 29 19:14:33 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need an explicit text format,
# otherwise Excel auto-converts the assigned string into a Number cell, which would
# not match the source data (these are meant to stay text, like the rest of column D).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D13", "D15", "D17", "D18", "D19", "D21", "D23", "D25", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D46", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '34.628.87'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '1.793.60'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '226.75'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = '0.557'
$ws.Range('E6').Value = '  +1.84%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '32.91'
$ws.Range('E8').Value = '  +3.19%  '
$ws.Range('D9').Value = '0.297'
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('D10').Value = '0.0694'
$ws.Range('E10').Value = '  +1.03%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '2.053.39'
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('D13').Value = '11.13'
$ws.Range('E13').Value = '  +0.90%  '
$ws.Range('D14').Value = '1.765.70'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').Value = '0.636'
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('D16').Value = '34.555.24'
$ws.Range('E16').Value = '  +1.26%  '
$ws.Range('D17').Value = '4.28'
$ws.Range('E17').Value = '  +2.42%  '
$ws.Range('D18').Value = '68.84'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').Value = '248.16'
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').Value = '0.0₃0802'
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('D21').Value = '11.26'
$ws.Range('E21').Value = '  +2.66%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '4.18'
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').Value = '165.05'
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('D27').Value = '16.57'
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('D28').Value = '0.116'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('E29').Value = '  +0.79%  '
$ws.Range('D30').Value = '4.14'
$ws.Range('E30').Value = '  +13.82%  '
$ws.Range('D31').Value = '3.83'
$ws.Range('E31').Value = '  +3.40%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0524'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.24'
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('D34').Value = '1.83'
$ws.Range('E34').Value = '  +1.87%  '
$ws.Range('D35').Value = '1.428.69'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('D36').Value = '2.60'
$ws.Range('E36').Value = '  +6.54%  '
$ws.Range('D37').Value = '0.671'
$ws.Range('E37').Value = '  +2.40%  '
$ws.Range('E38').Value = '  +2.12%  '
$ws.Range('D39').Value = '0.0192'
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('D40').Value = '85.27'
$ws.Range('E40').Value = '  +6.29%  '
$ws.Range('D41').Value = '2.40'
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '2.76'
$ws.Range('E42').Value = '  +2.81%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '0.932'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('E44').Value = '  +0.98%  '
$ws.Range('E45').Value = '  +3.67%  '
$ws.Range('D46').Value = '6.12'
$ws.Range('E46').Value = '  +0.89%  '
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('D48').Value = '1.953.78'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').Value = '106.21'
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('E51').Value = '  -4.53%  '
